$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2500888
$ws.Range("J17").Value = 2500888
$ws.Range("L17").Value = 7502664
$ws.Range("N17").Value = -7503000
$ws.Range("H112").Value = 1246.3962
$ws.Range("J112").Value = 1299.5625
$ws.Range("L112").Value = 3898.6875
$ws.Range("N112").Value = -6114.6875
$ws.Range("H116").Value = 4678.4443
$ws.Range("I116").Value = 6700
$ws.Range("J116").Value = 3061.2
$ws.Range("K116").Value = 6700
$ws.Range("L116").Value = 3061.2
$ws.Range("M116").Value = -3258
$ws.Range("N116").Value = -9945.200000000001
$ws.Range("H127").Value = 1808.5476
$ws.Range("I127").Value = 495.83334
$ws.Range("J127").Value = 2027.3334
$ws.Range("K127").Value = 1487.50002
$ws.Range("L127").Value = 6082.0002
$ws.Range("M127").Value = 3472.49998
$ws.Range("N127").Value = -16002.0002
$ws.Range("H129").Value = 737.3333
$ws.Range("I129").Value = 418.8
$ws.Range("J129").Value = 859.8461
$ws.Range("K129").Value = 1256.4
$ws.Range("L129").Value = 2579.5383
$ws.Range("M129").Value = 3743.6
$ws.Range("N129").Value = -12579.5383
$ws.Range("H138").Value = 1976.3684
$ws.Range("I138").Value = 1563.175
$ws.Range("J138").Value = 2948.5881
$ws.Range("K138").Value = 4689.525
$ws.Range("L138").Value = 8845.764299999999
$ws.Range("M138").Value = 450.4750000000004
$ws.Range("N138").Value = -19125.7643

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5609.383
$ws.Range("I32").Value = 4272.0435
$ws.Range("J32").Value = 13299.083
$ws.Range("K32").Value = 4272.0435
$ws.Range("L32").Value = 13299.083
$ws.Range("M32").Value = -3985.0435
$ws.Range("N32").Value = -13873.083
$ws.Range("H74").Value = 3274.6035
$ws.Range("I74").Value = 3530.5813
$ws.Range("J74").Value = 2540.8
$ws.Range("K74").Value = 3530.5813
$ws.Range("L74").Value = 2540.8
$ws.Range("M74").Value = -2656.5813
$ws.Range("N74").Value = -4288.8
$ws.Range("H77").Value = 3274.6035
$ws.Range("I77").Value = 3530.5813
$ws.Range("J77").Value = 2540.8
$ws.Range("K77").Value = 17652.9065
$ws.Range("L77").Value = 12704
$ws.Range("M77").Value = -13284.9065
$ws.Range("N77").Value = -21440
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4008
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 2091.303
$ws.Range("I132").Value = 1281.2727
$ws.Range("J132").Value = 3711.3635
$ws.Range("K132").Value = 3843.8181
$ws.Range("L132").Value = 11134.0905
$ws.Range("M132").Value = -1313.8181
$ws.Range("N132").Value = -16194.0905

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 24635.5
$ws.Range("J63").Value = 24635.5
$ws.Range("L63").Value = 24635.5
$ws.Range("N63").Value = -26007.5
$ws.Range("H66").Value = 24635.5
$ws.Range("J66").Value = 24635.5
$ws.Range("L66").Value = 73906.5
$ws.Range("N66").Value = -80770.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4578.375
$ws.Range("I31").Value = 2386.3809
$ws.Range("J31").Value = 5893.5713
$ws.Range("K31").Value = 2386.3809
$ws.Range("L31").Value = 5893.5713
$ws.Range("M31").Value = -2091.3809
$ws.Range("N31").Value = -6483.5713
$ws.Range("H34").Value = 4578.375
$ws.Range("I34").Value = 2386.3809
$ws.Range("J34").Value = 5893.5713
$ws.Range("K34").Value = 2386.3809
$ws.Range("L34").Value = 5893.5713
$ws.Range("M34").Value = -2184.3809
$ws.Range("N34").Value = -6297.5713
$ws.Range("H135").Value = 38911.61
$ws.Range("J135").Value = 38911.61
$ws.Range("L135").Value = 38911.61
$ws.Range("N135").Value = -49051.61

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 980
$ws.Range("I68").Value = 976
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 2928
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2117
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 980
$ws.Range("I71").Value = 976
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 8784
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -4728
$ws.Range("N71").Value = -17112
$ws.Range("H99").Value = 1798.625
$ws.Range("I99").Value = 1177.8
$ws.Range("J99").Value = 2833.3333
$ws.Range("K99").Value = 3533.4
$ws.Range("L99").Value = 8499.999899999999
$ws.Range("M99").Value = -1287.4
$ws.Range("N99").Value = -12991.9999
$ws.Range("H107").Value = 520225.44
$ws.Range("I107").Value = 526.8
$ws.Range("J107").Value = 845037.0600000001
$ws.Range("K107").Value = 1580.4
$ws.Range("L107").Value = 2535111.18
$ws.Range("M107").Value = 339.6000000000001
$ws.Range("N107").Value = -2538951.18

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 2009863.6
$ws.Range("I52").Value = 7338500
$ws.Range("J52").Value = 11625
$ws.Range("K52").Value = 7338500
$ws.Range("L52").Value = 11625
$ws.Range("M52").Value = -7338241
$ws.Range("N52").Value = -12143
$ws.Range("H102").Value = 2772.4333
$ws.Range("I102").Value = 2611.818
$ws.Range("J102").Value = 3214.125
$ws.Range("K102").Value = 2611.818
$ws.Range("L102").Value = 3214.125
$ws.Range("M102").Value = -989.8180000000002
$ws.Range("N102").Value = -6458.125
$ws.Range("H122").Value = 2359.1777
$ws.Range("I122").Value = 1768.6552
$ws.Range("J122").Value = 3429.5
$ws.Range("K122").Value = 5305.9656
$ws.Range("L122").Value = 10288.5
$ws.Range("M122").Value = -2855.9656
$ws.Range("N122").Value = -15188.5
$ws.Range("H126").Value = 2156.121
$ws.Range("I126").Value = 1861.7894
$ws.Range("J126").Value = 2555.5715
$ws.Range("K126").Value = 5585.3682
$ws.Range("L126").Value = 7666.7145
$ws.Range("M126").Value = -3115.3682
$ws.Range("N126").Value = -12606.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 48193.773
$ws.Range("I7").Value = 73818.14
$ws.Range("J7").Value = 3351.125
$ws.Range("K7").Value = 73818.14
$ws.Range("L7").Value = 3351.125
$ws.Range("M7").Value = -73706.14
$ws.Range("N7").Value = -3575.125
$ws.Range("H36").Value = 35000
$ws.Range("J36").Value = 35000
$ws.Range("L36").Value = 35000
$ws.Range("N36").Value = -36124
$ws.Range("H40").Value = 33696.438
$ws.Range("I40").Value = 47203.453
$ws.Range("K40").Value = 47203.453
$ws.Range("M40").Value = -47067.453
$ws.Range("H45").Value = 20500
$ws.Range("I45").Value = 8000
$ws.Range("K45").Value = 8000
$ws.Range("M45").Value = -7593
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H122").Value = 2235.9092
$ws.Range("I122").Value = 1801.2222
$ws.Range("J122").Value = 2536.8462
$ws.Range("K122").Value = 5403.6666
$ws.Range("L122").Value = 7610.5386
$ws.Range("M122").Value = -2953.6666
$ws.Range("N122").Value = -12510.5386
$ws.Range("H126").Value = 48193.773
$ws.Range("I126").Value = 73818.14
$ws.Range("J126").Value = 3351.125
$ws.Range("K126").Value = 221454.42
$ws.Range("L126").Value = 10053.375
$ws.Range("M126").Value = -218984.42
$ws.Range("N126").Value = -14993.375
$ws.Range("H136").Value = 27857.725
$ws.Range("I136").Value = 45922.74
$ws.Range("J136").Value = 3416.8235
$ws.Range("K136").Value = 137768.22
$ws.Range("L136").Value = 10250.4705
$ws.Range("M136").Value = -135218.22
$ws.Range("N136").Value = -15350.4705

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 870.8125
$ws.Range("I113").Value = 1010.2308
$ws.Range("J113").Value = 266.66666
$ws.Range("K113").Value = 3030.6924
$ws.Range("L113").Value = 799.9999799999999
$ws.Range("M113").Value = -860.6923999999999
$ws.Range("N113").Value = -5139.99998
$ws.Range("H122").Value = 47162.727
$ws.Range("I122").Value = 60357.707
$ws.Range("K122").Value = 181073.121
$ws.Range("M122").Value = -178623.121
$ws.Range("H126").Value = 112666.555
$ws.Range("I126").Value = 167382.33
$ws.Range("J126").Value = 3235
$ws.Range("K126").Value = 502146.99
$ws.Range("L126").Value = 9705
$ws.Range("M126").Value = -499676.99
$ws.Range("N126").Value = -14645
$ws.Range("H136").Value = 5468
$ws.Range("I136").Value = 952.2
$ws.Range("K136").Value = 2856.6
$ws.Range("M136").Value = -306.6000000000004

